# Generate Report for Handback
# Adds a new handback record (0f1b92a9-af8d-4718-a057-8fe11696b002) as row 4
# on the Overview, zh-cn and de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$newId   = "0f1b92a9-af8d-4718-a057-8fe11696b002"
$newHash = "a9047b220e0adce49fd493aaa884acd2c3e92473"

$mdName      = "$newId.md"
$zhXlfName   = "$newId.$newHash.zh-cn.xlf"
$deXlfName   = "$newId.$newHash.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$includeText  = "Include"

$zhHandoffDt  = "2016-02-23 08:44:28"
$zhHandbackDt = "2016-02-23 08:45:12"
$deHandoffDt  = "2016-02-23 08:44:39"
$deHandbackDt = "2016-02-23 08:45:32"

# -- External link targets (same host/repo naming convention used by the
#    other rows in this workbook). --
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/07c48a562fedd100bdeaed18bbf4696c1625ca6c/e2e/$mdName"
$zhHandoffUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47c1af899b6f34ec78acd3fcf953f9d156c8ca27/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName"
$zhMdForkUrl   = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f7f8c0e91c5db6f4d34ae5c862e751282f0d2abb/e2e/$mdName"
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a07daae279cc965c4fce4e987a7f6f7a68810d3f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName"
$deHandoffUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a25a6ff9151ffe508200280d1eb6149c0efb8ef8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName"
$deMdForkUrl   = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/57c700b53ce32577b9c51dd8e95d1093dd640d9a/e2e/$mdName"
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fa0094a9a7fc5e9965c0d4faa8c45cb450b5f7c6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName"

# ---------------------------------------------------------------------
# Sheet 1: Overview  (File Name | zh-cn | de-de)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsOverview.Range("A4").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B4").Value = $statusInSync
$wsZh.Range("D4").Value = $zhHandoffDt
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G4").Value = $zhHandbackDt
$wsZh.Range("H4").Value = $includeText

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsZh.Range("A4").Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("C4"), $zhHandoffUrl, [Type]::Missing, [Type]::Missing, $zhXlfName) | Out-Null
$wsZh.Range("C4").Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("E4"), $zhMdForkUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsZh.Range("E4").Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("F4"), $zhHandbackUrl, [Type]::Missing, [Type]::Missing, $zhXlfName) | Out-Null
$wsZh.Range("F4").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B4").Value = $statusInSync
$wsDe.Range("D4").Value = $deHandoffDt
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G4").Value = $deHandbackDt
$wsDe.Range("H4").Value = $includeText

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsDe.Range("A4").Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("C4"), $deHandoffUrl, [Type]::Missing, [Type]::Missing, $deXlfName) | Out-Null
$wsDe.Range("C4").Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("E4"), $deMdForkUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsDe.Range("E4").Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("F4"), $deHandbackUrl, [Type]::Missing, [Type]::Missing, $deXlfName) | Out-Null
$wsDe.Range("F4").Style = "HyperLink"

Write-Host "Handback row added for $newId"
